$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2756.3635
$ws.Range("I17").Value = 2727.7778
$ws.Range("J17").Value = 2885
$ws.Range("K17").Value = 8183.3334
$ws.Range("L17").Value = 8655
$ws.Range("M17").Value = -8015.3334
$ws.Range("N17").Value = -8991

$ws.Range("H28").Value = 1314.625
$ws.Range("I28").Value = 456.72726
$ws.Range("J28").Value = 3202
$ws.Range("K28").Value = 456.72726
$ws.Range("L28").Value = 3202
$ws.Range("M28").Value = 28.27274
$ws.Range("N28").Value = -4172

$ws.Range("H86").Value = 1366.6666
$ws.Range("I86").Value = 1366.6666
$ws.Range("K86").Value = 1366.6666
$ws.Range("M86").Value = -243.6666

$ws.Range("H89").Value = 1366.6666
$ws.Range("I89").Value = 1366.6666
$ws.Range("K89").Value = 6833.333000000001
$ws.Range("M89").Value = -1217.333000000001

$ws.Range("H98").Value = 1875.25
$ws.Range("I98").Value = 1875.25
$ws.Range("K98").Value = 1875.25
$ws.Range("M98").Value = -377.25

$ws.Range("H122").Value = 1875.25
$ws.Range("I122").Value = 1875.25
$ws.Range("K122").Value = 5625.75
$ws.Range("M122").Value = -3175.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1999.5
$ws.Range("I45").Value = 1999
$ws.Range("J45").Value = 2000
$ws.Range("K45").Value = 1999
$ws.Range("L45").Value = 2000
$ws.Range("M45").Value = -1622
$ws.Range("N45").Value = -2754

$ws.Range("H110").Value = 62500976
$ws.Range("I110").Value = 1205.5
$ws.Range("J110").Value = 125000750
$ws.Range("K110").Value = 1205.5
$ws.Range("L110").Value = 125000750
$ws.Range("M110").Value = 839.5
$ws.Range("N110").Value = -125004840

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 843.7778
$ws.Range("I22").Value = 849.3333
$ws.Range("J22").Value = 832.6667
$ws.Range("K22").Value = 849.3333
$ws.Range("L22").Value = 832.6667
$ws.Range("M22").Value = -676.3333
$ws.Range("N22").Value = -1178.6667

$ws.Range("H42").Value = 199999
$ws.Range("J42").Value = 199999
$ws.Range("L42").Value = 199999
$ws.Range("N42").Value = -200655

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 590446.8
$ws.Range("I99").Value = 558040.9
$ws.Range("J99").Value = 626903.5
$ws.Range("K99").Value = 558040.9
$ws.Range("L99").Value = 626903.5
$ws.Range("M99").Value = -556542.9
$ws.Range("N99").Value = -629899.5

$ws.Range("H122").Value = 1285.5264
$ws.Range("I122").Value = 1402.0667
$ws.Range("J122").Value = 848.5
$ws.Range("K122").Value = 4206.2001
$ws.Range("L122").Value = 2545.5
$ws.Range("M122").Value = -1756.2001
$ws.Range("N122").Value = -7445.5

$ws.Range("H126").Value = 590446.8
$ws.Range("I126").Value = 558040.9
$ws.Range("J126").Value = 626903.5
$ws.Range("K126").Value = 1674122.7
$ws.Range("L126").Value = 1880710.5
$ws.Range("M126").Value = -1671652.7
$ws.Range("N126").Value = -1885650.5

$ws.Range("H141").Value = 685943.4399999999
$ws.Range("J141").Value = 685943.4399999999
$ws.Range("L141").Value = 685943.4399999999
$ws.Range("N141").Value = -696303.4399999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H59").Value = 30
$ws.Range("I59").Value = 20
$ws.Range("J59").Value = 40
$ws.Range("K59").Value = 60
$ws.Range("L59").Value = 120
$ws.Range("M59").Value = 480
$ws.Range("N59").Value = -1200

$ws.Range("H107").Value = 853.8570999999999
$ws.Range("I107").Value = 418.25
$ws.Range("K107").Value = 1254.75
$ws.Range("M107").Value = 665.25

$ws.Range("H131").Value = 2131.2856
$ws.Range("I131").Value = 979.75
$ws.Range("J131").Value = 3666.6667
$ws.Range("K131").Value = 2939.25
$ws.Range("L131").Value = 11000.0001
$ws.Range("M131").Value = 2100.75
$ws.Range("N131").Value = -21080.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1056.5714
$ws.Range("I102").Value = 897
$ws.Range("J102").Value = 2014
$ws.Range("K102").Value = 897
$ws.Range("L102").Value = 2014
$ws.Range("M102").Value = 725
$ws.Range("N102").Value = -5258

$ws.Range("H122").Value = 1106.5217
$ws.Range("I122").Value = 1092.8572
$ws.Range("K122").Value = 3278.5716
$ws.Range("M122").Value = -828.5715999999998

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1190.8
$ws.Range("I7").Value = 996
$ws.Range("J7").Value = 1483
$ws.Range("K7").Value = 996
$ws.Range("L7").Value = 1483
$ws.Range("M7").Value = -884
$ws.Range("N7").Value = -1707

$ws.Range("H22").Value = 258.33334
$ws.Range("J22").Value = 300
$ws.Range("L22").Value = 300
$ws.Range("N22").Value = -890

$ws.Range("H27").Value = 258.33334
$ws.Range("J27").Value = 300
$ws.Range("L27").Value = 300
$ws.Range("N27").Value = -514

$ws.Range("H34").Value = 14500
$ws.Range("I34").Value = 14500
$ws.Range("K34").Value = 14500
$ws.Range("M34").Value = -14328

$ws.Range("H40").Value = 3166.3333
$ws.Range("I40").Value = 2999.5
$ws.Range("K40").Value = 2999.5
$ws.Range("M40").Value = -2863.5

$ws.Range("H46").Value = 2888
$ws.Range("I46").Value = 1198.4
$ws.Range("K46").Value = 1198.4
$ws.Range("M46").Value = -1010.4

$ws.Range("H55").Value = 1350.8235
$ws.Range("I55").Value = 1161.6666
$ws.Range("J55").Value = 1804.8
$ws.Range("K55").Value = 1161.6666
$ws.Range("L55").Value = 1804.8
$ws.Range("M55").Value = -988.6666
$ws.Range("N55").Value = -2150.8

$ws.Range("H61").Value = 2001.3334
$ws.Range("I61").Value = 1502
$ws.Range("K61").Value = 1502
$ws.Range("M61").Value = -1300

$ws.Range("H113").Value = 2001.3334
$ws.Range("I113").Value = 1502
$ws.Range("K113").Value = 1502
$ws.Range("M113").Value = 668

$ws.Range("H122").Value = 5272.375
$ws.Range("I122").Value = 4106.75
$ws.Range("K122").Value = 12320.25
$ws.Range("M122").Value = -9870.25

$ws.Range("H126").Value = 1190.8
$ws.Range("I126").Value = 996
$ws.Range("J126").Value = 1483
$ws.Range("K126").Value = 2988
$ws.Range("L126").Value = 4449
$ws.Range("M126").Value = -518
$ws.Range("N126").Value = -9389

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 503.2857
$ws.Range("I107").Value = 504
$ws.Range("J107").Value = 499
$ws.Range("K107").Value = 1512
$ws.Range("L107").Value = 1497
$ws.Range("M107").Value = 408
$ws.Range("N107").Value = -5337

$ws.Range("H122").Value = 1627.4
$ws.Range("I122").Value = 1621.2222
$ws.Range("J122").Value = 1636.6666
$ws.Range("K122").Value = 4863.6666
$ws.Range("L122").Value = 4909.9998
$ws.Range("M122").Value = -2413.6666
$ws.Range("N122").Value = -9809.9998
